# Add a new "2022-Q3" quarter sheet right after the "总计" (total) sheet,
# and update the 总计 summary table with the new quarter's row.

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item(1)
$prevQuarterSheet = $wb.Worksheets.Item(2)   # currently "2022-Q2" — same column layout as the new sheet

# ---------------------------------------------------------------------
# 1) Duplicate the existing "2022-Q2" sheet (to inherit its exact column
#    layout / styles) right after "总计", rename it to "2022-Q3", then
#    overwrite its contents with the 2022-Q3 fund data.
# ---------------------------------------------------------------------
$prevQuarterSheet.Copy([System.Type]::Missing, $totalSheet)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# The old sheet had 28 data rows (rows 2-29); 2022-Q3 only has 23
# (rows 2-24), so drop the extra tail rows.
$newSheet.Range("A25:H29").Clear()

$codes    = @("008283","011152","008515","011468","870009","010365","011914","011913","501025","012170","001703","872019","014768","011469","006810","006809","014767","014364","012289","011647","012288","012171","011648")
$names    = @("易方达金融行业股票","富兰克林国海兴海回报混合","富兰克林国海基本面优选混合","国富竞争优势三年持有期混合A","广发资管平衡精选一年持有混合A","鹏华港股通中证香港银行投资指数（LOF）C","华夏永泓一年持有混合C","华夏永泓一年持有混合A","鹏华港股通中证香港银行投资指数（LOF）A","华夏永顺一年持有混合A","银华沪港深增长股票A","广发资管平衡精选一年持有混合C","景顺华城稳健6月持有混合C","国富竞争优势三年持有期混合C","泰康港股通中证香港银行投资指数C","泰康港股通中证香港银行投资指数A","景顺华城稳健6月持有混合A","银华沪港深增长股票C","泰康沪港深成长混合C","博时港股通红利精选混合A","泰康沪港深成长混合A","华夏永顺一年持有混合C","博时港股通红利精选混合C")
$sizes    = @("17.50","14.87","11.82","11.65","7.47","4.92","11.19","11.18","2.38","8.93","2.07","1.09","2.25","0.67","0.73","0.69","1.26","0.25","0.44","0.11","0.20","0.38","0.02")
$stockPos = @("86.30","84.99","85.30","82.69","92.14","94.13","31.67","31.67","94.13","28.10","88.44","92.14","24.78","82.69","94.14","94.14","24.78","88.44","72.20","82.44","72.20","28.10","82.44")
$ratio    = @("8.25","5.01","5.25","3.70","5.11","3.38","0.93","0.93","3.38","0.79","3.09","5.11","1.79","3.70","3.39","3.39","1.79","3.09","1.63","5.60","1.63","0.79","5.60")
$value    = @("1.4438","0.7450","0.6206","0.4310","0.3817","0.1663","0.1041","0.1040","0.0804","0.0705","0.0640","0.0557","0.0403","0.0248","0.0247","0.0234","0.0226","0.0077","0.0072","0.0062","0.0033","0.0030","0.0011")
$rank     = @(5,5,4,6,8,10,8,8,10,10,9,8,3,6,10,10,3,9,6,5,6,10,5)

for ($i = 0; $i -lt $codes.Count; $i++) {
    $r = $i + 2
    $newSheet.Cells.Item($r, 1).Value = $i

    $newSheet.Cells.Item($r, 2).NumberFormat = "@"
    $newSheet.Cells.Item($r, 2).Value = $codes[$i]

    $newSheet.Cells.Item($r, 3).Value = $names[$i]

    $newSheet.Cells.Item($r, 4).NumberFormat = "@"
    $newSheet.Cells.Item($r, 4).Value = $sizes[$i]

    $newSheet.Cells.Item($r, 5).NumberFormat = "@"
    $newSheet.Cells.Item($r, 5).Value = $stockPos[$i]

    $newSheet.Cells.Item($r, 6).NumberFormat = "@"
    $newSheet.Cells.Item($r, 6).Value = $ratio[$i]

    $newSheet.Cells.Item($r, 7).NumberFormat = "@"
    $newSheet.Cells.Item($r, 7).Value = $value[$i]

    $newSheet.Cells.Item($r, 8).Value = $rank[$i]
}

# ---------------------------------------------------------------------
# 2) Update the "总计" (total) summary sheet: insert a new row for the
#    2022-Q3 quarter at the top of the data (row 2), pushing the rest
#    down, and renumber the A-column index.
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

# Reuse row 3's (the row that just got pushed down, still carrying the
# original formatting) cell format for the newly-inserted row 2 cells.
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A2:D2").PasteSpecial(-4122)  # xlPasteFormats

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 23
$totalSheet.Cells.Item(2, 4).Value = 4.43

# Renumber the index column (A) for the rows that got pushed down.
for ($r = 3; $r -le 9; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

# Restore original view state (总计 sheet selected, cell A1 active).
$totalSheet.Activate() | Out-Null
$totalSheet.Range("A1").Select() | Out-Null
